$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-04-14 Sunday" "2024-04-15 Monday"

Replace-Text "23×40=920" "75×97=7275"
Replace-Text "30×33=990" "93×79=7347"
Replace-Text "98×83=8134" "61×80=4880"
Replace-Text "22×28=616" "38×38=1444"
Replace-Text "91×91=8281" "63×41=2583"
Replace-Text "14×72=1008" "94×92=8648"
Replace-Text "38×62=2356" "84×69=5796"
Replace-Text "76×65=4940" "99×58=5742"
Replace-Text "87×19=1653" "68×46=3128"
Replace-Text "49×47=2303" "62×19=1178"
Replace-Text "72×88=6336" "85×34=2890"
Replace-Text "40×20=800" "27×75=2025"
Replace-Text "26×30=780" "43×90=3870"
Replace-Text "48×38=1824" "76×98=7448"
Replace-Text "15×14=210" "31×75=2325"
Replace-Text "17×82=1394" "96×73=7008"
Replace-Text "89×64=5696" "68×30=2040"
Replace-Text "18×43=774" "34×41=1394"
Replace-Text "41×13=533" "86×15=1290"
Replace-Text "40×90=3600" "17×56=952"
Replace-Text "97×69=6693" "12×17=204"
Replace-Text "73×22=1606" "73×65=4745"
Replace-Text "43×81=3483" "31×95=2945"
Replace-Text "33×51=1683" "99×70=6930"
Replace-Text "36×48=1728" "41×75=3075"
